# Refresh the crypto price/volume snapshot (GitHub Actions bot run).
# Only the cells whose Price (D) / Volume(1h) (E) figures moved, or whose
# row got re-ranked (swapping Coin name / Link / Price / Volume with the
# neighboring row), are touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '34.168.92'
$ws.Cells.Item(2, 5).Value = '  +0.91%  '
$ws.Cells.Item(3, 4).Value = '1.780.75'
$ws.Cells.Item(3, 5).Value = '  +0.06%  '
$ws.Cells.Item(4, 5).Value = '  +0.12%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '225.99'
$ws.Cells.Item(5, 5).Value = '  +0.72%  '
$ws.Cells.Item(6, 5).Value = '  +0.08%  '
$ws.Cells.Item(7, 5).Value = '  +0.14%  '
$ws.Cells.Item(8, 5).Value = '  +0.07%  '
$ws.Cells.Item(9, 5).Value = '  +1.12%  '
$ws.Cells.Item(10, 5).Value = '  +2.30%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0945'
$ws.Cells.Item(11, 5).Value = '  +1.02%  '
$ws.Cells.Item(12, 4).Value = '2.038.23'
$ws.Cells.Item(12, 5).Value = '  +0.47%  '
$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(13, 4).Value = '1.798.93'
$ws.Cells.Item(13, 5).Value = '  +1.45%  '
$ws.Cells.Item(14, 2).Value = 'Chainlink'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '10.93'
$ws.Cells.Item(14, 5).Value = '  -1.99%  '
$ws.Cells.Item(15, 2).Value = 'Polygon'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.624'
$ws.Cells.Item(15, 5).Value = '  +2.46%  '
$ws.Cells.Item(16, 2).Value = 'WrappedBTC'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(16, 4).Value = '34.153.50'
$ws.Cells.Item(16, 5).Value = '  +0.87%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '4.19'
$ws.Cells.Item(17, 5).Value = '  +1.28%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '67.90'
$ws.Cells.Item(18, 5).Value = '  +1.69%  '
$ws.Cells.Item(19, 4).Value = '0.0₃0803'
$ws.Cells.Item(19, 5).Value = '  +4.24%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '246.01'
$ws.Cells.Item(20, 5).Value = '  +3.04%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '10.98'
$ws.Cells.Item(21, 5).Value = '  +4.14%  '
$ws.Cells.Item(22, 5).Value = '  +0.12%  '
$ws.Cells.Item(23, 5).Value = '  +1.99%  '
$ws.Cells.Item(24, 5).Value = '  -1.00%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '162.45'
$ws.Cells.Item(25, 5).Value = '  +1.02%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '7.18'
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '16.29'
$ws.Cells.Item(27, 5).Value = '  +1.39%  '
$ws.Cells.Item(28, 5).Value = '  +1.97%  '
$ws.Cells.Item(29, 5).Value = '  +0.33%  '
$ws.Cells.Item(30, 2).Value = 'Hedera'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '0.0521'
$ws.Cells.Item(30, 5).Value = '  +2.46%  '
$ws.Cells.Item(31, 2).Value = 'PancakeSwap'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '1.23'
$ws.Cells.Item(31, 5).Value = '  +0.83%  '
$ws.Cells.Item(32, 5).Value = '  +4.25%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '3.71'
$ws.Cells.Item(33, 5).Value = '  +5.21%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '1.80'
$ws.Cells.Item(34, 5).Value = '  -0.97%  '
$ws.Cells.Item(35, 4).Value = '1.441.31'
$ws.Cells.Item(35, 5).Value = '  +3.81%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.658'
$ws.Cells.Item(36, 5).Value = '  +3.92%  '
$ws.Cells.Item(38, 5).Value = '  +2.89%  '
$ws.Cells.Item(39, 5).Value = '  +0.38%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '80.24'
$ws.Cells.Item(40, 5).Value = '  +2.60%  '
$ws.Cells.Item(41, 5).Value = '  -0.11%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.924'
$ws.Cells.Item(42, 5).Value = '  +1.41%  '
$ws.Cells.Item(43, 5).Value = '  +0.92%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '13.51'
$ws.Cells.Item(44, 5).Value = '  +0.45%  '
$ws.Cells.Item(45, 2).Value = 'Kaspa'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.0510'
$ws.Cells.Item(45, 5).Value = '  +0.18%  '
$ws.Cells.Item(46, 2).Value = 'FraxShare'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '6.08'
$ws.Cells.Item(46, 5).Value = '  +3.91%  '
$ws.Cells.Item(47, 5).Value = '  -0.01%  '
$ws.Cells.Item(48, 5).Value = '  -2.89%  '
$ws.Cells.Item(49, 4).Value = '1.939.76'
$ws.Cells.Item(49, 5).Value = '  +0.17%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '104.46'
$ws.Cells.Item(50, 5).Value = '  -1.06%  '
$ws.Cells.Item(51, 5).Value = '  +0.12%  '
